# Electricity sector calibration: set geothermal, biomass, and natural gas
# steam turbine (CHP) to guaranteed dispatch (BAU Guaranteed Dispatch = 1
# for all years), and document why on the About sheet.

$wb  = $excel.ActiveWorkbook
$wsAbout   = $wb.Worksheets.Item("About")
$wsBgdpbes = $wb.Worksheets.Item("BGDPbES")

# --- About sheet: add explanatory note -------------------------------------
$wsAbout.Range("A28").Value = "Natural gas steam turbine (representing CHP), biomass, and geothermal all run at fixed capacity factors, so"
$wsAbout.Range("A29").Value = "we flag them here."

# --- BGDPbES sheet: guarantee full dispatch (value 1) for the years 2015-2050
# Row 3  = natural gas steam turbine (representing CHP)
# Row 10 = biomass
# Row 11 = geothermal
$wsBgdpbes.Range("B3:AK3").Value = 1
$wsBgdpbes.Range("B10:AK10").Value = 1
$wsBgdpbes.Range("B11:AK11").Value = 1

# --- Restore the selections left behind by the edits (matches where the
# author's cursor ended up after typing the new rows/values) -----------------
$wsBgdpbes.Range("B10:AK11").Select() | Out-Null
$wsAbout.Range("A30").Select() | Out-Null
